$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new traded row (row 8) ---
$ws.Range("A8").Value = 42654.743761574071
$ws.Range("A8").NumberFormat = "m/d/yy h:mm"
$ws.Range("B8").Value = $False
$ws.Range("C8").Value = 9811.24
$ws.Range("D8").Value = 9818.6
$ws.Range("E8").Value = 104.43
$ws.Range("F8").Value = 104.269997
$ws.Range("G8").Value = $False
$ws.Range("H8").Value = -0.15
$ws.Range("I8").Value = $False

# --- Re-fit the column widths now that the repeater added a new, differently
#     sized row of data underneath the existing trades (bestFit columns resize
#     themselves in Excel when new data is entered) ---
$ws.Columns.Item(1).ColumnWidth = 14.5
$ws.Columns.Item(2).ColumnWidth = 7.333333333333333
$ws.Columns.Item(3).ColumnWidth = 7.0
$ws.Columns.Item(4).ColumnWidth = 10.333333333333334
$ws.Columns.Item(5).ColumnWidth = 10.0
$ws.Columns.Item(6).ColumnWidth = 10.0
$ws.Columns.Item(7).ColumnWidth = 9.5
$ws.Columns.Item(8).ColumnWidth = 13.833333333333334
$ws.Columns.Item(9).ColumnWidth = 11.0
